# Weekly update: insert 3 new "Palta" price records (Hass, all 3 qualities)
# for Agrícola del Norte S.A. de Arica, ahead of the existing historical
# rows (which all shift down by 3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows right before the current row 243, pushing the
# existing 243:270 data block down to 246:273.
$ws.Rows("243:245").Insert()

# Shared, constant attributes for every "Palta" record on this sheet.
$mercadoId = 1
$mercado   = "Agrícola del Norte S.A. de Arica"
$region    = "Arica y Parinacota"
$codreg    = 15
$tipo      = "Fruta"
$productoId = 100106
$producto  = "Oleaginosos"
$categoriaId = 100106002
$categoria = "Palta"
$unidad    = "`$/bandeja 10 kilos"
$origen    = "Perú"
$kgUnidad  = 10

# --- New row 243: Palta / Hass / Primera ---
$ws.Range("A243").Value = $mercadoId
$ws.Range("B243").Value = $mercado
$ws.Range("C243").Value = $region
$ws.Range("D243").Value = 45166
$ws.Range("E243").Value = $codreg
$ws.Range("F243").Value = $tipo
$ws.Range("G243").Value = $productoId
$ws.Range("H243").Value = $producto
$ws.Range("I243").Value = $categoriaId
$ws.Range("J243").Value = $categoria
$ws.Range("K243").Value = "Hass"
$ws.Range("L243").Value = "Primera"
$ws.Range("M243").Value = 380
$ws.Range("N243").Value = 28000
$ws.Range("O243").Value = 29000
$ws.Range("P243").Value = 28474
$ws.Range("Q243").Value = $unidad
$ws.Range("R243").Value = $origen
$ws.Range("S243").Value = 2847
$ws.Range("T243").Value = $kgUnidad

# --- New row 244: Palta / Hass / Segunda ---
$ws.Range("A244").Value = $mercadoId
$ws.Range("B244").Value = $mercado
$ws.Range("C244").Value = $region
$ws.Range("D244").Value = 45166
$ws.Range("E244").Value = $codreg
$ws.Range("F244").Value = $tipo
$ws.Range("G244").Value = $productoId
$ws.Range("H244").Value = $producto
$ws.Range("I244").Value = $categoriaId
$ws.Range("J244").Value = $categoria
$ws.Range("K244").Value = "Hass"
$ws.Range("L244").Value = "Segunda"
$ws.Range("M244").Value = 410
$ws.Range("N244").Value = 26000
$ws.Range("O244").Value = 27000
$ws.Range("P244").Value = 26463
$ws.Range("Q244").Value = $unidad
$ws.Range("R244").Value = $origen
$ws.Range("S244").Value = 2646
$ws.Range("T244").Value = $kgUnidad

# --- New row 245: Palta / Hass / Tercera ---
$ws.Range("A245").Value = $mercadoId
$ws.Range("B245").Value = $mercado
$ws.Range("C245").Value = $region
$ws.Range("D245").Value = 45166
$ws.Range("E245").Value = $codreg
$ws.Range("F245").Value = $tipo
$ws.Range("G245").Value = $productoId
$ws.Range("H245").Value = $producto
$ws.Range("I245").Value = $categoriaId
$ws.Range("J245").Value = $categoria
$ws.Range("K245").Value = "Hass"
$ws.Range("L245").Value = "Tercera"
$ws.Range("M245").Value = 350
$ws.Range("N245").Value = 24000
$ws.Range("O245").Value = 25000
$ws.Range("P245").Value = 24571
$ws.Range("Q245").Value = $unidad
$ws.Range("R245").Value = $origen
$ws.Range("S245").Value = 2457
$ws.Range("T245").Value = $kgUnidad
